$wb = $excel.ActiveWorkbook

# --- Metadata sheet (sheet 1): insert a "Jurisdiction" row after "Contact" ---
$ws = $wb.Worksheets.Item(1)

# Shift existing rows 11-14 (Description, Purpose, Copyright, Immutable) down to 12-15,
# working from the bottom up so we don't clobber values before they are copied.
for ($r = 14; $r -ge 11; $r--) {
    $dst = $r + 1
    $ws.Cells.Item($dst, 1).Value = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($dst, 2).Value = $ws.Cells.Item($r, 2).Value()
}

# New row 11: Jurisdiction / (empty)
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Row 15 now holds what used to be row 14's content; carry over its formatting too
# (value copies above did not bring formatting along for the newly-extended row).
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

# --- Update the Date value on the Metadata sheet ---
$ws.Cells.Item(8, 2).Value = "2024-09-17T19:55:11+00:00"

# --- Rename the second sheet ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Include #0"
